# Rename the inline picture shapes (logo images) that live in the
# document's header/footer parts:
#   - footer1.xml (first-page footer, docPr id=3) : PearsonLogo  image2.png -> image1.png
#   - footer2.xml (primary footer,   docPr id=2) : PearsonLogo  image2.png -> image1.png
#   - header1.xml (first-page header, docPr id=1) : BTec_Logo-Orange image1.jpg -> image2.jpg
#
# Each rename is a structural edit, so the Section/HeaderFooter handle is
# re-fetched fresh from ActiveDocument before every access below instead of
# being reused across edits.

$d = $word.ActiveDocument

# --- First-page header: BTec_Logo-Orange ---------------------------------
$hdrFirst = $d.Sections.Item(1).Headers.Item(2)   # wdHeaderFooterFirstPage
if ($hdrFirst.Range.InlineShapes.Count -ge 1) {
    $btecLogo = $hdrFirst.Range.InlineShapes.Item(1)
    $btecLogo.Name = "image2.jpg"
}

# --- Primary footer: PearsonLogo -----------------------------------------
$ftrPrimary = $d.Sections.Item(1).Footers.Item(1)  # wdHeaderFooterPrimary
if ($ftrPrimary.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoPrimary = $ftrPrimary.Range.InlineShapes.Item(1)
    $pearsonLogoPrimary.Name = "image1.png"
}

# --- First-page footer: PearsonLogo ---------------------------------------
$ftrFirst = $d.Sections.Item(1).Footers.Item(2)   # wdHeaderFooterFirstPage
if ($ftrFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoFirst = $ftrFirst.Range.InlineShapes.Item(1)
    $pearsonLogoFirst.Name = "image1.png"
}
